# Generate Report for Handoff
# - Flip the per-language status from "Handed back: in sync with en-US" to
#   "Ready for handoff" on the Overview sheet and on each language sheet.
# - Bump the associated "generated/handoff" timestamps forward a few seconds
#   (new xliff generation pass).
# - The Status/zh-cn/de-de columns got noticeably narrower now that the new
#   status text is shorter than the old one (re-fit column width).

$wb = $excel.ActiveWorkbook

# ----- Overview sheet -----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-16 04:56:32"

$ov.Columns.Item(5).ColumnWidth = 16.3
$ov.Columns.Item(6).ColumnWidth = 16.3

# ----- zh-cn sheet -----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("H2").Value = "2016-08-16 04:56:27"

$zh.Columns.Item(3).ColumnWidth = 16.3

# ----- de-de sheet -----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Ready for handoff"
$de.Range("H2").Value = "2016-08-16 04:56:32"

$de.Columns.Item(3).ColumnWidth = 16.3
